$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Portfolio")

# --- Header row tweaks ---------------------------------------------------
# "Value date" column is dropped; "P&L" shifts from K1 into J1, and a new
# "Buy/sell" header takes K1 (with a centered, boxed style like the rest of
# row 1, but only left/right borders).
$ws.Range("J1").Value = "P&L"
$ws.Range("K1").Value = "Buy/sell"

# Give K1 the same look as the other header cells (centered, bold white on
# navy, like the existing header style) but trim the border to left/right
# only.
$ws.Range("K1").Copy()
$ws.Range("K1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("K1").Borders.Item(8).LineStyle = -4142  # xlEdgeTop -> none
$ws.Range("K1").Borders.Item(9).LineStyle = -4142  # xlEdgeBottom -> none

# --- Row 2: "Z" replaces the old numeric Qty value in H2, J2 becomes a
# plain number (no more date format), K2 gets a Buy/sell flag ------------
$ws.Range("H2").Value = "Z"
$ws.Range("J2").Value = 5
$ws.Range("J2").ClearFormats()
$ws.Range("K2").Value = 1

# --- Buy/sell flags for the remaining rows -------------------------------
$ws.Range("K3").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("K5").Value = 1
$ws.Range("K6").Value = 0
$ws.Range("K7").Value = 1

# --- New "P&L" worksheet, placed after Portfolio -------------------------
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$newSheet.Name = "P&L"

$ws.Range("H2").Select() | Out-Null
